$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 229
$ws.Range("D229").Value = 44516
$ws.Range("J229").Value = 750
$ws.Range("K229").Value = 12000
$ws.Range("L229").Value = 12000
$ws.Range("M229").Value = 12000
$ws.Range("O229").Value = 'Región de Ñuble'
$ws.Range("P229").Value = 600

# Row 230
$ws.Range("D230").Value = 44270
$ws.Range("J230").Value = 150
$ws.Range("K230").Value = 7000
$ws.Range("L230").Value = 7000
$ws.Range("M230").Value = 7000
$ws.Range("P230").Value = 350

# Row 231
$ws.Range("K231").Value = 6000
$ws.Range("L231").Value = 6000
$ws.Range("M231").Value = 6000
$ws.Range("O231").Value = 'Provincia de Llanquihue'
$ws.Range("P231").Value = 300

# Row 232
$ws.Range("D232").Value = 44295
$ws.Range("J232").Value = 350
$ws.Range("K232").Value = 8000
$ws.Range("L232").Value = 8000
$ws.Range("M232").Value = 8000
$ws.Range("P232").Value = 400

# Row 233
$ws.Range("D233").Value = 44217
$ws.Range("J233").Value = 200
$ws.Range("K233").Value = 9000
$ws.Range("L233").Value = 10000
$ws.Range("M233").Value = 9500
$ws.Range("P233").Value = 475

# Row 234
$ws.Range("D234").Value = 44509
$ws.Range("J234").Value = 700
$ws.Range("K234").Value = 12000
$ws.Range("L234").Value = 12500
$ws.Range("M234").Value = 12250
$ws.Range("O234").Value = 'Región de Ñuble'
$ws.Range("P234").Value = 612

# Row 235
$ws.Range("J235").Value = 350
$ws.Range("K235").Value = 6500
$ws.Range("L235").Value = 6500
$ws.Range("M235").Value = 6500
$ws.Range("O235").Value = 'Provincia de Llanquihue'
$ws.Range("P235").Value = 325

# Row 236
$ws.Range("D236").Value = 44421
$ws.Range("J236").Value = 400
$ws.Range("K236").Value = 8000
$ws.Range("L236").Value = 8000
$ws.Range("M236").Value = 8000
$ws.Range("O236").Value = 'Región de Ñuble'
$ws.Range("P236").Value = 400

# Row 237
$ws.Range("K237").Value = 6000
$ws.Range("L237").Value = 6000
$ws.Range("M237").Value = 6000
$ws.Range("O237").Value = 'Provincia de Llanquihue'
$ws.Range("P237").Value = 300

# Row 238
$ws.Range("D238").Value = 44383
$ws.Range("J238").Value = 350
$ws.Range("K238").Value = 8000
$ws.Range("L238").Value = 8000
$ws.Range("M238").Value = 8000
$ws.Range("O238").Value = 'Región de Ñuble'
$ws.Range("P238").Value = 400

# Row 239
$ws.Range("D239").Value = 44244
$ws.Range("J239").Value = 27
$ws.Range("K239").Value = 7500
$ws.Range("L239").Value = 7500
$ws.Range("M239").Value = 7500
$ws.Range("P239").Value = 375

# Row 240
$ws.Range("D240").Value = 44307
$ws.Range("J240").Value = 70
$ws.Range("K240").Value = 6000
$ws.Range("L240").Value = 6000
$ws.Range("M240").Value = 6000
$ws.Range("P240").Value = 300

# Row 241
$ws.Range("D241").Value = 44273
$ws.Range("K241").Value = 7000
$ws.Range("L241").Value = 7000
$ws.Range("M241").Value = 7000
$ws.Range("P241").Value = 350

# Row 242
$ws.Range("D242").Value = 44433
$ws.Range("J242").Value = 150
$ws.Range("K242").Value = 6500
$ws.Range("L242").Value = 6500
$ws.Range("M242").Value = 6500
$ws.Range("P242").Value = 325

# Row 243
$ws.Range("D243").Value = 44302
$ws.Range("J243").Value = 750
$ws.Range("K243").Value = 6000
$ws.Range("L243").Value = 6000
$ws.Range("M243").Value = 6000
$ws.Range("O243").Value = 'Provincia de Llanquihue'
$ws.Range("P243").Value = 300

# Row 244
$ws.Range("D244").Value = 44179
$ws.Range("J244").Value = 150
$ws.Range("K244").Value = 8000
$ws.Range("L244").Value = 8000
$ws.Range("M244").Value = 8000
$ws.Range("O244").Value = 'Región de Coquimbo'
$ws.Range("P244").Value = 400

# Row 245
$ws.Range("A245").Value = 4
$ws.Range("B245").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C245").Value = 'Los Lagos'
$ws.Range("D245").Value = 44491
$ws.Range("D245").NumberFormat = $ws.Range("D244").NumberFormat()
$ws.Range("E245").Value = 10
$ws.Range("F245").Value = 100114013
$ws.Range("G245").Value = 'Zanahoria'
$ws.Range("H245").Value = 'Sin especificar'
$ws.Range("I245").Value = 'Primera'
$ws.Range("J245").Value = 750
$ws.Range("K245").Value = 12000
$ws.Range("L245").Value = 12000
$ws.Range("M245").Value = 12000
$ws.Range("N245").Value = '$/saco 20 kilos'
$ws.Range("O245").Value = 'Región de Ñuble'
$ws.Range("P245").Value = 600
$ws.Range("Q245").Value = 20
$ws.Range("R245").Value = 'Hortaliza'
